$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update employee names in column A (rows 2-24), reflecting reordering of shared strings
$ws.Range("A2").Value = "Aline Castro"
$ws.Range("A3").Value = "Tamires Teixeira"
$ws.Range("A4").Value = "Daniela Fernandes"
$ws.Range("A5").Value = "Adriana Hunhoff"
$ws.Range("A6").Value = "Natalia Farias"
$ws.Range("A7").Value = "Luana Umpierre"
$ws.Range("A8").Value = "Nicolas Silva"
$ws.Range("A9").Value = "Vania Fagundes"
$ws.Range("A10").Value = "Carlla Bo"
$ws.Range("A11").Value = "Julio Acauan"
$ws.Range("A12").Value = "Josue Alos"
$ws.Range("A13").Value = "Alexia Pereira"
$ws.Range("A14").Value = "Daniel Machado"
$ws.Range("A15").Value = "Brenda Pereira"
$ws.Range("A16").Value = "Amanda Bernardes"
$ws.Range("A17").Value = "Dominique Daudt"
$ws.Range("A18").Value = "Eduarda Santos"
$ws.Range("A19").Value = "Michele Mattidorff"
$ws.Range("A20").Value = "Jonathan Cardoso"
$ws.Range("A21").Value = "Gabriel Wolff"
$ws.Range("A22").Value = "Gabriel Winter"
$ws.Range("A23").Value = "Igor Martins"
$ws.Range("A24").Value = "Brenda Fossa"

# Update numeric data for rows 2-13 and totals row 25 (columns B:J)
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 85
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 101

$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 64
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 16
$ws.Range("G3").Value = 13
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 98

$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 15
$ws.Range("E4").Value = 24
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 63

$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 48
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 61

$ws.Range("B6").Value = 19
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 22
$ws.Range("E6").Value = 19
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 51

$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 21
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = 13
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 43

$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 12
$ws.Range("F8").Value = 29
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 43

$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 19
$ws.Range("F9").Value = 17
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 5
$ws.Range("J9").Value = 41

$ws.Range("B10").Value = 4
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 11
$ws.Range("F10").Value = 23
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 34

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 16

$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 2

$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0

$ws.Range("B25").Value = 71
$ws.Range("C25").Value = 38
$ws.Range("D25").Value = 234
$ws.Range("E25").Value = 128
$ws.Range("F25").Value = 131
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = 2
$ws.Range("I25").Value = 15
$ws.Range("J25").Value = 553

